$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing voter code for row 4 (registered voter who hasn't voted)
$ws.Range("A4").Value = 444444

# Move selection to D10 to match the saved cursor position
$ws.Range("D10").Select()
